$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Bump the term package version
$ws.Range("B3").Value = "1.1.0"

# Update the published date to match the new version
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"
